# Daten aktualisiert am 2024-04-05
# Append 5 new ticker rows to the bottom of the existing list on Sheet1
# (A464:A468), growing the used range from A1:A463 to A1:A468.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTickers = @("IMX-USD", "MNT-USD", "PEPE-USD", "GRT-USD", "TAO-USD")

$startRow = 464
for ($i = 0; $i -lt $newTickers.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newTickers[$i]
}
